# Relationales Modell.xlsx - add Sebastian's "Bauteile/Roboter" table set
# as a second worksheet, rename the original sheet to "Markus".
#
# Commit message:
#   Einfuegen der Tabellen durch Sebastian Golchert
#   Anlegen eines Issues fuer Fragen zu den Datentypen

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the existing sheet "Tabelle1" -> "Markus"
# ------------------------------------------------------------------
$markus = $wb.Worksheets.Item(1)
$markus.Name = "Markus"

# ------------------------------------------------------------------
# 2. Add the new sheet "Sebastian" right after "Markus"
# ------------------------------------------------------------------
$sebastian = $wb.Worksheets.Add($null, $markus)
$sebastian.Name = "Sebastian"

# Column A is noticeably wider on this sheet (longest label:
# "Roboterkomponenten")
$sebastian.Columns.Item(1).ColumnWidth = 18.166666666666668

# Match the 2 cm top/bottom page margins used on the rest of the workbook
$sebastian.PageSetup.TopMargin = 56.692913399999995
$sebastian.PageSetup.BottomMargin = 56.692913399999995

# ------------------------------------------------------------------
# 3. Header row (same header used on "Markus")
# ------------------------------------------------------------------
$sebastian.Cells.Item(1,1).Value = "Spalte"
$sebastian.Cells.Item(1,2).Value = "Datentyp `n[Länge]"
$sebastian.Cells.Item(1,2).WrapText = $true
$sebastian.Cells.Item(1,3).Value = "Null-`nOption"
$sebastian.Cells.Item(1,3).WrapText = $true
$sebastian.Cells.Item(1,5).Value = "Constraints"
$sebastian.Cells.Item(1,6).Value = "Bemerkungen"

# ------------------------------------------------------------------
# 4. "Bauteile" table
# ------------------------------------------------------------------
$sebastian.Cells.Item(3,1).Value = "Bauteile"
$sebastian.Cells.Item(3,1).Font.Bold = $true
$sebastian.Cells.Item(3,1).Font.Color = 255

$sebastian.Cells.Item(4,1).Value = "TeileNr"
$sebastian.Cells.Item(4,1).Font.Bold = $true
$sebastian.Cells.Item(4,2).Value = "int"
$sebastian.Cells.Item(4,3).Value = "not null"
$sebastian.Cells.Item(4,4).Value = "auto_increment"
$sebastian.Cells.Item(4,5).Value = "Primary Key"

$sebastian.Cells.Item(5,1).Value = "Bezeichnung"
$sebastian.Cells.Item(5,2).Value = "varchar[?]"
$sebastian.Cells.Item(5,3).Value = "not null"

# ------------------------------------------------------------------
# 5. "Roboterkomponenten" table
# ------------------------------------------------------------------
$sebastian.Cells.Item(8,1).Value = "Roboterkomponenten"
$sebastian.Cells.Item(8,1).Font.Bold = $true
$sebastian.Cells.Item(8,1).Font.Color = 255

$sebastian.Cells.Item(9,1).Value = "Stückzahl"
$sebastian.Cells.Item(9,2).Value = "int"
$sebastian.Cells.Item(9,3).Value = "not null"

# ------------------------------------------------------------------
# 6. "Roboter" table
# ------------------------------------------------------------------
$sebastian.Cells.Item(13,1).Value = "Roboter"
$sebastian.Cells.Item(13,1).Font.Bold = $true
$sebastian.Cells.Item(13,1).Font.Color = 255

$sebastian.Cells.Item(14,1).Value = "RoboterID"
$sebastian.Cells.Item(14,1).Font.Bold = $true
$sebastian.Cells.Item(14,2).Value = "int"
$sebastian.Cells.Item(14,3).Value = "not null"
$sebastian.Cells.Item(14,4).Value = "auto_increment"
$sebastian.Cells.Item(14,5).Value = "Primary Key"
$sebastian.Cells.Item(14,6).Value = "?"

$sebastian.Cells.Item(15,1).Value = "Bezeichnung"
$sebastian.Cells.Item(15,2).Value = "varchar[?]"
$sebastian.Cells.Item(15,3).Value = "not null"

$sebastian.Cells.Item(16,1).Value = "Materialkosten"
$sebastian.Cells.Item(16,2).Value = "decimal(7,2)"
$sebastian.Cells.Item(16,3).Value = "not null"

$sebastian.Cells.Item(17,1).Value = "Produktionskosten"
$sebastian.Cells.Item(17,2).Value = "decimal(7,2)"
$sebastian.Cells.Item(17,3).Value = "not null"

# ------------------------------------------------------------------
# 7. Selections: "Markus" keeps a leftover range selection from the
#    copy/paste that seeded "Sebastian"; "Sebastian" ends up with the
#    cursor on the last-typed cell and becomes the active tab.
# ------------------------------------------------------------------
$markus.Range("A1:F15").Select()
$sebastian.Range("C17").Select()
$sebastian.Activate()
